# enable_disable_notification_settings.xlsx
# "add api test case first..." — while drafting the new API test case the
# author resized row 5 (the SYMENADISNOT-002 test case row) so its
# multi-line "Steps" text is fully visible, and had scrolled further down
# the frozen-header list (top visible row moves from A44 to A46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow row 5 so the wrapped text is fully visible (15.75pt -> 50.4pt).
$ws.Rows.Item(5).RowHeight = 50.4

# Scroll the frozen (bottom-left) pane so row 46 is the first visible row.
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
